$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1399.69
$summary.Range("B4").Value = -0.32
$summary.Range("B5").Value = -0.46
$summary.Range("B6").Value = 14
$summary.Range("B7").Value = 5
$summary.Range("B9").Value = 35.71

# --- Sheet: Strategy Status ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 99.69
$status.Range("D5").Value = 14
$status.Range("E5").Value = -0.32
$status.Range("F5").Value = -0.31
$status.Range("G5").Value = 35.71

# --- New trade row to append to "All Trades" and "MarketMaking" sheets ---
# Note: column B holds a date-looking string ("2026-02-17") that must stay
# plain text (matching the rest of the column), so it is written with a
# leading apostrophe to stop Excel from auto-converting it to a date serial.
$newRow = @(14, "'2026-02-17", "20:03:26", "MarketMaking", "UP", 0.84, 0.85, "CLOSED", 1.1905, 0.01, 99.69, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.13)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($i = 0; $i -lt $newRow.Length; $i++) {
        $ws.Cells.Item(15, $i + 1).Value = $newRow[$i]
    }
}
